$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 259 - this shifts the existing rows
# 259..306 down to 260..307 (matching the old data verbatim), and the
# dimension grows to A1:R307 automatically.
$ws.Rows(259).Insert()

# Populate the freshly inserted row 259 with a new price record, copying
# the unchanged columns (A,B,C,E,F,G,H,I,O,R) from the row that used to
# occupy this slot (now at row 260) and setting the new reported values.
$ws.Cells.Item(259, 1).Value2 = 9
$ws.Cells.Item(259, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(259, 3).Value2 = "Metropolitana"
$ws.Cells.Item(259, 4).Value2 = 44816
$ws.Cells.Item(259, 5).Value2 = 13
$ws.Cells.Item(259, 6).Value2 = 100112001
$ws.Cells.Item(259, 7).Value2 = "Berenjena"
$ws.Cells.Item(259, 8).Value2 = "Sin especificar"
$ws.Cells.Item(259, 9).Value2 = "Primera"
$ws.Cells.Item(259, 10).Value2 = 160
$ws.Cells.Item(259, 11).Value2 = 14000
$ws.Cells.Item(259, 12).Value2 = 14000
$ws.Cells.Item(259, 13).Value2 = 14000
$ws.Cells.Item(259, 14).Value2 = '$/caja 50 unidades'
$ws.Cells.Item(259, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(259, 16).Value2 = 280
$ws.Cells.Item(259, 17).Value2 = 50
$ws.Cells.Item(259, 18).Value2 = "Hortaliza"
